$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab from "ShearF-HW50.xpc" to "ShearF"
$ws.Name = "ShearF"

# Add new row 16 mirroring row 15's pattern (formats + values)
$ws.Range("A15:M15").Copy()
$ws.Range("A16:M16").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = $ws.Range("B15").Text
$ws.Range("C16:M16").Value = 1

"done"
